$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$roster = @(
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Dillon Brooks", "SG,SF", "Houston Rockets"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Bobby Portis", "PF,C", "Milwaukee Bucks"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
    @("Bilal Coulibaly", "SG,SF", "Washington Wizards"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Dennis Schröder", "PG,SG", "Golden State Warriors"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Jonathan Kuminga", "SF,PF", "Golden State Warriors"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Paul George", "SF,PF", "Philadelphia 76ers")
)

for ($i = 0; $i -lt $roster.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $roster[$i][0]
    $ws.Cells.Item($row, 2).Value = $roster[$i][1]
    $ws.Cells.Item($row, 3).Value = $roster[$i][2]
}
